$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.569.03"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.433.02"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.19"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.25"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.595"
$ws.Range("E7").Value = "  -2.48%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.694"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  +7.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.06"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.85"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.441.42"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.472.26"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.56"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("E19").Value = "  +11.60%  "
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "84.54"
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "312.26"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.82"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.75"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.10"
$ws.Range("E27").Value = "  -5.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.77"
$ws.Range("E28").Value = "  +4.07%  "
$ws.Range("E29").Value = "  +4.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "44.78"
$ws.Range("E30").Value = "  +4.98%  "
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.74"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.323"
$ws.Range("E39").Value = "  +12.61%  "
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "142.89"
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.82"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.37"
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.107.23"
$ws.Range("E48").Value = "  -2.00%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.97"
$ws.Range("E49").Value = "  +3.71%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.30"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.10"
$ws.Range("E51").Value = "  +29.69%  "
